$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.816.11"
$ws.Range("E2").Value = "  -0.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.038.18"
$ws.Range("E3").Value = "  -0.85%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.40"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.61"
$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0821"
$ws.Range("E10").Value = "  +2.08%  "

$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.71"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.340.22"
$ws.Range("E13").Value = "  -0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.06"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.778"
$ws.Range("E15").Value = "  +2.76%  "

$ws.Range("E16").Value = "  -2.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.048.76"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.788.15"
$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("E19").Value = "  -0.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.89"
$ws.Range("E20").Value = "  -6.57%  "

$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.91"
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -1.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  +2.78%  "

$ws.Range("E26").Value = "  +1.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.99"
$ws.Range("E27").Value = "  +1.39%  "

$ws.Range("E28").Value = "  -3.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.82"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("E30").Value = "  -2.18%  "

$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("E32").Value = "  +8.81%  "

$ws.Range("E33").Value = "  -3.15%  "

$ws.Range("E34").Value = "  -1.65%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.48"
$ws.Range("E36").Value = "  +0.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  +2.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.42"
$ws.Range("E38").Value = "  +4.44%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.02"
$ws.Range("E40").Value = "  +5.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.535.00"
$ws.Range("E41").Value = "  +1.05%  "

$ws.Range("E42").Value = "  +0.72%  "

$ws.Range("E43").Value = "  -1.47%  "

$ws.Range("E44").Value = "  -2.06%  "

$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("E46").Value = "  -2.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.06"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("E48").Value = "  -0.45%  "

$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.03"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.229.09"
$ws.Range("E51").Value = "  -0.85%  "
